# chore: (normalization) added support for managers table in the frontend
#
# All rows (2-6) on the active sheet get their "Reporting Manager" column
# (column C) normalized to a single manager name, "Rajesh".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C6").Value = "Rajesh"

# Update the active selection to C6, matching the recorded sheet view state.
$ws.Range("C6").Select()
